# Update workbook for data through 2021-12-04 (adds November 26 data point)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab / sheet name
$ws.Name = "Through 2021-11-26"

# Update the row label for November
$ws.Range("A12").Value = "November (through 11-26)"

# Update November row (row 12) values
$ws.Range("B12").Value = 28
$ws.Range("C12").Value = 66
$ws.Range("D12").Value = 98
$ws.Range("E12").Value = 58
$ws.Range("F12").Value = 46
$ws.Range("G12").Value = 182
$ws.Range("H12").Value = 177

# Update Total row (row 13) values
$ws.Range("B13").Value = 286
$ws.Range("C13").Value = 552
$ws.Range("D13").Value = 808
$ws.Range("E13").Value = 673
$ws.Range("F13").Value = 528
$ws.Range("G13").Value = 1239
$ws.Range("H13").Value = 1620
